$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (10th column); this shifts J:T to K:U.
$ws.Columns("J").Insert()

# The hyperlink on the old J2 (now carried by the engine still at J2) needs
# to move to its new home at K2 - drop the stale one first (deleting via a
# range-scoped Hyperlinks collection clears the sheet's hyperlinks) then
# re-add it pointing at the cell that now holds the FTP host value.
$ws.Range("J2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("K2"), "ftp://ftp.jenne.com/")

# Update the lookup column label (was "ManPartNum", now "ItemNum").
$ws.Range("C2").Value = "ItemNum"

# New header + value for the freshly inserted column.
$ws.Range("J1").Value = "Source"
$ws.Range("J2").Value = "Local"

# Match the saved selection state from the diff.
$ws.Range("J3").Select()
